# Insert a new price-record row for "Red Globe" grapes on 2022-02-08
# (serial 44606) at row 54 of the weekly ("semanal") Uva sheet for
# "Terminal La Palmera de La Serena". All rows that were previously at
# 54..84 shift down to 55..85 (their data is untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 54..84 down to 55..85, opening up a blank row 54.
$ws.Rows.Item(54).Insert()

# Populate the newly opened row 54 with the new record.
$ws.Cells.Item(54, 1).Value  = 8
$ws.Cells.Item(54, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(54, 3).Value  = "Coquimbo"
$ws.Cells.Item(54, 4).Value  = 44606
$ws.Cells.Item(54, 5).Value  = 4
$ws.Cells.Item(54, 6).Value  = "Fruta"
$ws.Cells.Item(54, 7).Value  = 100109
$ws.Cells.Item(54, 8).Value  = "Uva"
$ws.Cells.Item(54, 9).Value  = 100109001
$ws.Cells.Item(54, 10).Value = "Uva"
$ws.Cells.Item(54, 11).Value = "Red Globe"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 400
$ws.Cells.Item(54, 14).Value = 11500
$ws.Cells.Item(54, 15).Value = 12000
$ws.Cells.Item(54, 16).Value = 11750
$ws.Cells.Item(54, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(54, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(54, 19).Value = 653
$ws.Cells.Item(54, 20).Value = 18

# Match the date-formatted style already used by the rest of column D.
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
